$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "967×4="
$t.Cell(1,2).Range.Text = "794×9="
$t.Cell(1,3).Range.Text = "674×2="
$t.Cell(1,4).Range.Text = "206×7="
$t.Cell(1,5).Range.Text = "543×8="
$t.Cell(5,1).Range.Text = "452×6="
$t.Cell(5,2).Range.Text = "508×8="
$t.Cell(5,3).Range.Text = "171×7="
$t.Cell(5,4).Range.Text = "453×3="
$t.Cell(5,5).Range.Text = "725×9="
$t.Cell(10,1).Range.Text = "977×3="
$t.Cell(10,2).Range.Text = "367×4="
$t.Cell(10,3).Range.Text = "598×7="
$t.Cell(10,4).Range.Text = "623×5="
$t.Cell(10,5).Range.Text = "610×5="
$t.Cell(15,1).Range.Text = "262×2="
$t.Cell(15,2).Range.Text = "619×8="
$t.Cell(15,3).Range.Text = "122×8="
$t.Cell(15,4).Range.Text = "982×3="
$t.Cell(15,5).Range.Text = "350×8="
$t.Cell(20,1).Range.Text = "678×3="
$t.Cell(20,2).Range.Text = "756×2="
$t.Cell(20,3).Range.Text = "333×4="
$t.Cell(20,4).Range.Text = "562×5="
$t.Cell(20,5).Range.Text = "349×3="

Write-Output "done"
